$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cyclic rotation of stock-card batch data (columns B, D, E, F, G)
# across grouped rows that represent the same product, per the target revision.

# Group rows [146, 147, 148]
$ws.Range("B146").Value = 64350
$ws.Range("D146").Value = 66.44
$ws.Range("E146").Value = 70.63
$ws.Range("F146").Value = 2
$ws.Range("G146").Value = 132.88
$ws.Range("B147").Value = 57756
$ws.Range("D147").Value = 66.44
$ws.Range("E147").Value = 79.37
$ws.Range("F147").Value = -100
$ws.Range("G147").Value = -6644
$ws.Range("B148").Value = 53925
$ws.Range("D148").Value = 66.44
$ws.Range("E148").Value = 79.37
$ws.Range("F148").Value = 1
$ws.Range("G148").Value = 66.44

# Group rows [233, 234]
$ws.Range("B233").Value = 64979
$ws.Range("D233").Value = 295.75
$ws.Range("E233").Value = 314.41
$ws.Range("F233").Value = 0
$ws.Range("G233").Value = 0
$ws.Range("B234").Value = 48719
$ws.Range("D234").Value = 295.75
$ws.Range("E234").Value = 353.35
$ws.Range("F234").Value = -81
$ws.Range("G234").Value = -23955.75

# Group rows [246, 247]
$ws.Range("B246").Value = 48706
$ws.Range("D246").Value = 33.3
$ws.Range("E246").Value = 39.8
$ws.Range("F246").Value = -144
$ws.Range("G246").Value = -4795.2
$ws.Range("B247").Value = 64973
$ws.Range("D247").Value = 33.3
$ws.Range("E247").Value = 35.4
$ws.Range("F247").Value = 53
$ws.Range("G247").Value = 1764.9

# Group rows [292, 293]
$ws.Range("B292").Value = 55373
$ws.Range("D292").Value = 144.28
$ws.Range("E292").Value = 163.62
$ws.Range("F292").Value = -94
$ws.Range("G292").Value = -13562.32
$ws.Range("B293").Value = 63520
$ws.Range("D293").Value = 144.28
$ws.Range("E293").Value = 153.4
$ws.Range("F293").Value = 73
$ws.Range("G293").Value = 10532.44

# Group rows [294, 295]
$ws.Range("B294").Value = 63571
$ws.Range("D294").Value = 143.48
$ws.Range("E294").Value = 152.53
$ws.Range("F294").Value = 4
$ws.Range("G294").Value = 573.92
$ws.Range("B295").Value = 57802
$ws.Range("D295").Value = 143.48
$ws.Range("E295").Value = 162.71
$ws.Range("F295").Value = -79
$ws.Range("G295").Value = -11334.92

# Group rows [299, 300]
$ws.Range("B299").Value = 55356
$ws.Range("D299").Value = 47.64
$ws.Range("E299").Value = 54.04
$ws.Range("F299").Value = -158
$ws.Range("G299").Value = -7527.12
$ws.Range("B300").Value = 63510
$ws.Range("D300").Value = 47.64
$ws.Range("E300").Value = 50.66
$ws.Range("F300").Value = 135
$ws.Range("G300").Value = 6431.4

# Group rows [311, 312]
$ws.Range("B311").Value = 61605
$ws.Range("D311").Value = 111.96
$ws.Range("E311").Value = 133.78
$ws.Range("F311").Value = -13
$ws.Range("G311").Value = -1455.48
$ws.Range("B312").Value = 63563
$ws.Range("D312").Value = 111.96
$ws.Range("E312").Value = 119.04
$ws.Range("F312").Value = 2
$ws.Range("G312").Value = 223.92

# Group rows [420, 421]
$ws.Range("B420").Value = 58047
$ws.Range("D420").Value = 105.54
$ws.Range("E420").Value = 126.1
$ws.Range("F420").Value = 42
$ws.Range("G420").Value = 4432.68
$ws.Range("B421").Value = 47097
$ws.Range("D421").Value = 112.28
$ws.Range("E421").Value = 134.16
$ws.Range("F421").Value = 15
$ws.Range("G421").Value = 1684.2

# Group rows [465, 466]
$ws.Range("B465").Value = 65069
$ws.Range("D465").Value = 13.45
$ws.Range("E465").Value = 14.3
$ws.Range("F465").Value = 2
$ws.Range("G465").Value = 26.9
$ws.Range("B466").Value = 53757
$ws.Range("D466").Value = 13.45
$ws.Range("E466").Value = 16.08
$ws.Range("F466").Value = -159
$ws.Range("G466").Value = -2138.55

# Group rows [472, 473]
$ws.Range("B472").Value = 45695
$ws.Range("D472").Value = 19.73
$ws.Range("E472").Value = 23.58
$ws.Range("F472").Value = -36
$ws.Range("G472").Value = -710.28
$ws.Range("B473").Value = 64915
$ws.Range("D473").Value = 19.73
$ws.Range("E473").Value = 20.98
$ws.Range("F473").Value = 0
$ws.Range("G473").Value = 0

# Group rows [485, 486]
$ws.Range("B485").Value = 45709
$ws.Range("D485").Value = 13.15
$ws.Range("E485").Value = 15.69
$ws.Range("F485").Value = -300
$ws.Range("G485").Value = -3945
$ws.Range("B486").Value = 64925
$ws.Range("D486").Value = 13.15
$ws.Range("E486").Value = 13.97
$ws.Range("F486").Value = 179
$ws.Range("G486").Value = 2353.85

# Group rows [490, 491]
$ws.Range("B490").Value = 65067
$ws.Range("D490").Value = 14.73
$ws.Range("E490").Value = 15.65
$ws.Range("F490").Value = 235
$ws.Range("G490").Value = 3461.55
$ws.Range("B491").Value = 53595
$ws.Range("D491").Value = 14.73
$ws.Range("E491").Value = 17.61
$ws.Range("F491").Value = -335
$ws.Range("G491").Value = -4934.55

# Group rows [570, 571]
$ws.Range("B570").Value = 64810
$ws.Range("D570").Value = 273.92
$ws.Range("E570").Value = 291.22
$ws.Range("F570").Value = 5
$ws.Range("G570").Value = 1369.6
$ws.Range("B571").Value = 53319
$ws.Range("D571").Value = 273.92
$ws.Range("E571").Value = 310.64
$ws.Range("F571").Value = -6
$ws.Range("G571").Value = -1643.52

# Group rows [602, 603]
$ws.Range("B602").Value = 64830
$ws.Range("D602").Value = 32.83
$ws.Range("E602").Value = 34.9
$ws.Range("F602").Value = 111
$ws.Range("G602").Value = 3644.13
$ws.Range("B603").Value = 60022
$ws.Range("D603").Value = 32.83
$ws.Range("E603").Value = 37.22
$ws.Range("F603").Value = -113
$ws.Range("G603").Value = -3709.79

# Group rows [711, 712]
$ws.Range("B711").Value = 63150
$ws.Range("D711").Value = 75.68000000000001
$ws.Range("E711").Value = 80.45
$ws.Range("F711").Value = 51
$ws.Range("G711").Value = 3859.68
$ws.Range("B712").Value = 61428
$ws.Range("D712").Value = 69.16
$ws.Range("E712").Value = 73.52
$ws.Range("F712").Value = 1
$ws.Range("G712").Value = 69.16
